$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Resolving-Mac" sending-cluster rows (original rows 4 and 5).
# Deleting row 4 shifts the old row 5 up into row 4, so delete at index 4 twice.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Refresh values for the remaining rows (2 and 3) with the new TPM-derived numbers.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("O2").Value = 0.06816352613805679
$ws.Range("P2").Value = 0.0681635261380568
$ws.Range("S2").Value = 0.06816352613805679
$ws.Range("T2").Value = 0.0681635261380568

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("M3").Value = 1.873282666666666
$ws.Range("N3").Value = 5.619847999999999
$ws.Range("O3").Value = 0.9318364738619431
$ws.Range("P3").Value = 0.9318364738619432
$ws.Range("Q3").Value = 0.3106870524044444
$ws.Range("R3").Value = 2.79618347164
$ws.Range("S3").Value = 0.9318364738619431
$ws.Range("T3").Value = 0.9318364738619432

Write-Output "done"
